$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text interpretation for numeric-looking values (e.g. "1.000"),
# then restore the default "Normal" style so cell formatting matches the original.
$textRange = $ws.Range("D2:E51")
$textRange.NumberFormat = "@"

$ws.Range("D2").Value = '30.291.76'
$ws.Range("E2").Value = '  -0.16%  '
$ws.Range("D3").Value = '1.929.45'
$ws.Range("E3").Value = '  -0.19%  '
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").Value = '0.7559'
$ws.Range("E5").Value = '  +5.45%  '
$ws.Range("D6").Value = '244.45'
$ws.Range("E6").Value = '  -2.61%  '
$ws.Range("D7").Value = '0.9995'
$ws.Range("E7").Value = '  -0.16%  '
$ws.Range("D8").Value = '0.3178'
$ws.Range("E8").Value = '  -2.98%  '
$ws.Range("D9").Value = '27.48'
$ws.Range("E9").Value = '  -0.42%  '
$ws.Range("D10").Value = '0.06985'
$ws.Range("E10").Value = '  -2.77%  '
$ws.Range("D11").Value = '0.7800'
$ws.Range("E11").Value = '  -2.76%  '
$ws.Range("D12").Value = '0.07981'
$ws.Range("E12").Value = '  -1.19%  '
$ws.Range("D13").Value = '1.930.26'
$ws.Range("E13").Value = '  -0.10%  '
$ws.Range("D14").Value = '5.350'
$ws.Range("E14").Value = '  -1.29%  '
$ws.Range("D15").Value = '94.24'
$ws.Range("E15").Value = '  -0.34%  '
$ws.Range("D16").Value = '14.40'
$ws.Range("E16").Value = '  -3.40%  '
$ws.Range("D17").Value = '30.290.24'
$ws.Range("E17").Value = '  -0.15%  '
$ws.Range("D18").Value = '251.91'
$ws.Range("E18").Value = '  -0.28%  '
$ws.Range("D19").Value = '0.000007910'
$ws.Range("E19").Value = '  -2.90%  '
$ws.Range("E20").Value = '  -1.44%  '
$ws.Range("D21").Value = '2.185.02'
$ws.Range("E21").Value = '  +0.02%  '
$ws.Range("D22").Value = '0.9990'
$ws.Range("E22").Value = '  -0.23%  '
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("D24").Value = '6.669'
$ws.Range("E24").Value = '  -3.89%  '
$ws.Range("D25").Value = '9.478'
$ws.Range("E25").Value = '  -2.62%  '
$ws.Range("D26").Value = '165.79'
$ws.Range("E26").Value = '  +0.33%  '
$ws.Range("B27").Value = 'Stellar'
$ws.Range("C27").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D27").Value = '0.1330'
$ws.Range("E27").Value = '  +2.97%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '18.90'
$ws.Range("E28").Value = '  -1.72%  '
$ws.Range("D29").Value = '2.207'
$ws.Range("E29").Value = '  -5.50%  '
$ws.Range("D30").Value = '1.365'
$ws.Range("E30").Value = '  +0.09%  '
$ws.Range("E31").Value = '  -2.03%  '
$ws.Range("D32").Value = '4.367'
$ws.Range("E32").Value = '  -1.26%  '
$ws.Range("D33").Value = '4.110'
$ws.Range("E33").Value = '  -2.21%  '
$ws.Range("D34").Value = '0.05156'
$ws.Range("E34").Value = '  -1.03%  '
$ws.Range("D35").Value = '1.276'
$ws.Range("E35").Value = '  +0.78%  '
$ws.Range("D36").Value = '0.7449'
$ws.Range("E36").Value = '  -0.31%  '
$ws.Range("D37").Value = '2.769'
$ws.Range("E37").Value = '  +0.21%  '
$ws.Range("D38").Value = '0.01945'
$ws.Range("E38").Value = '  -0.83%  '
$ws.Range("D39").Value = '2.793'
$ws.Range("E39").Value = '  -0.16%  '
$ws.Range("D40").Value = '77.68'
$ws.Range("E40").Value = '  -1.68%  '
$ws.Range("D41").Value = '6.404'
$ws.Range("E41").Value = '  -1.07%  '
$ws.Range("D42").Value = '0.4461'
$ws.Range("E42").Value = '  -1.47%  '
$ws.Range("D43").Value = '1.961'
$ws.Range("E43").Value = '  -3.13%  '
$ws.Range("D44").Value = '0.9994'
$ws.Range("E44").Value = '  -0.15%  '
$ws.Range("D45").Value = '0.8318'
$ws.Range("E45").Value = '  -0.94%  '
$ws.Range("D46").Value = '100.79'
$ws.Range("E46").Value = '  -1.08%  '
$ws.Range("D47").Value = '9.717'
$ws.Range("E47").Value = '  -0.74%  '
$ws.Range("D48").Value = '7.439'
$ws.Range("E48").Value = '  +0.34%  '
$ws.Range("D49").Value = '983.42'
$ws.Range("E49").Value = '  +11.41%  '
$ws.Range("D50").Value = '37.22'
$ws.Range("E50").Value = '  +1.52%  '
$ws.Range("D51").Value = '0.06005'
$ws.Range("E51").Value = '  -0.96%  '

$textRange.Style = "Normal"
